$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 18:31:25"
$ws.Range("A3").Value = "Total filas: 306"
$ws.Range("A23").Value = "06:17:28"
$ws.Range("C23").Value = "16_SANTA ANA"
$ws.Range("D23").Value = 64
$ws.Range("A24").Value = "05:57:13"
$ws.Range("C24").Value = "23_HERNANDEZ"
$ws.Range("D24").Value = 84
$ws.Range("C40").Value = "11_ETCHEVERRY"
$ws.Range("C41").Value = "15_ABASTO"
$ws.Range("A111").Value = "11:54:18"
$ws.Range("C111").Value = "15X38_ABASTO"
$ws.Range("D111").Value = 0
$ws.Range("A112").Value = "10:50:41"
$ws.Range("C112").Value = "23_HERNANDEZ"
$ws.Range("D112").Value = 64
$ws.Range("A113").Value = "11:52:01"
$ws.Range("C113").Value = "225_GOMEZ"
$ws.Range("D113").Value = 2
$ws.Range("C118").Value = "15_ABASTO"
$ws.Range("C119").Value = "16_P MOR-SANTA ANA"
$ws.Range("C120").Value = "15_ABASTO"
$ws.Range("C121").Value = "16_P MOR-SANTA ANA"
$ws.Range("A130").Value = "11:34:59"
$ws.Range("C130").Value = "15_ABASTO"
$ws.Range("D130").Value = 59
$ws.Range("A131").Value = "11:47:17"
$ws.Range("C131").Value = "14_ABASTO"
$ws.Range("D131").Value = 46
$ws.Range("C139").Value = "14_ABASTO"
$ws.Range("C141").Value = "15X38_ABASTO"
$ws.Range("A158").Value = "11:47:17"
$ws.Range("C158").Value = "16_P MOR-SANTA ANA"
$ws.Range("D158").Value = 98
$ws.Range("A159").Value = "12:11:52"
$ws.Range("C159").Value = "23_HERNANDEZ"
$ws.Range("D159").Value = 74
$ws.Range("A167").Value = "12:11:52"
$ws.Range("C167").Value = "10_OLMOS"
$ws.Range("D167").Value = 110
$ws.Range("A168").Value = "12:45:56"
$ws.Range("C168").Value = "23_HERNANDEZ"
$ws.Range("D168").Value = 76
$ws.Range("A216").Value = "14:12:26"
$ws.Range("C216").Value = "14_ABASTO"
$ws.Range("D216").Value = 113
$ws.Range("A217").Value = "15:17:33"
$ws.Range("C217").Value = "16_SANTA ANA"
$ws.Range("D217").Value = 48
$ws.Range("A244").Value = "16:14:44"
$ws.Range("C244").Value = "215C_EL PATO"
$ws.Range("D244").Value = 55
$ws.Range("A245").Value = "15:46:07"
$ws.Range("C245").Value = "23_HERNANDEZ"
$ws.Range("D245").Value = 83
$ws.Range("C273").Value = "16_P MOR-SANTA ANA"
$ws.Range("C274").Value = "15_ABASTO"
$ws.Range("A282").Value = "18:31:25"
$ws.Range("B282").Value = "18:34"
$ws.Range("C282").Value = "14_ABASTO"
$ws.Range("D282").Value = 3
$ws.Range("A284").Value = "16:37:06"
$ws.Range("B284").Value = "18:36"
$ws.Range("C284").Value = "15X38_ABASTO"
$ws.Range("D284").Value = 119
$ws.Range("A285").Value = "17:36:10"
$ws.Range("B285").Value = "18:37"
$ws.Range("C285").Value = "23_HERNANDEZ"
$ws.Range("D285").Value = 61
$ws.Range("A286").Value = "18:31:25"
$ws.Range("B286").Value = "18:40"
$ws.Range("C286").Value = "23_HERNANDEZ"
$ws.Range("D286").Value = 9
$ws.Range("B287").Value = "18:41"
$ws.Range("C287").Value = "10_OLMOS"
$ws.Range("D287").Value = 88
$ws.Range("A288").Value = "16:52:42"
$ws.Range("B288").Value = "18:45"
$ws.Range("C288").Value = "16_SANTA ANA"
$ws.Range("D288").Value = 113
$ws.Range("B289").Value = "18:52"
$ws.Range("C289").Value = "17_ROMERO"
$ws.Range("D289").Value = 99
$ws.Range("A290").Value = "17:13:39"
$ws.Range("B290").Value = "18:57"
$ws.Range("C290").Value = "16_P MOR-SANTA ANA"
$ws.Range("D290").Value = 104
$ws.Range("B291").Value = "18:59"
$ws.Range("C291").Value = "14_ABASTO"
$ws.Range("D291").Value = 106
$ws.Range("B292").Value = "19:00"
$ws.Range("C292").Value = "14_ABASTO"
$ws.Range("D292").Value = 84
$ws.Range("A293").Value = "17:13:39"
$ws.Range("B293").Value = "19:03"
$ws.Range("C293").Value = "215_EL PELIGRO"
$ws.Range("D293").Value = 110
$ws.Range("A294").Value = "17:36:10"
$ws.Range("B294").Value = "19:04"
$ws.Range("C294").Value = "215_EL PELIGRO"
$ws.Range("D294").Value = 88
$ws.Range("A295").Value = "17:56:03"
$ws.Range("B295").Value = "19:10"
$ws.Range("D295").Value = 74
$ws.Range("A296").Value = "18:12:30"
$ws.Range("B296").Value = "19:10"
$ws.Range("C296").Value = "16_SANTA ANA"
$ws.Range("D296").Value = 58
$ws.Range("A297").Value = "17:48:33"
$ws.Range("B297").Value = "19:12"
$ws.Range("D297").Value = 84
$ws.Range("A298").Value = "17:56:03"
$ws.Range("B298").Value = "19:16"
$ws.Range("C298").Value = "17_ROMERO"
$ws.Range("D298").Value = 80
$ws.Range("A299").Value = "18:12:30"
$ws.Range("B299").Value = "19:16"
$ws.Range("C299").Value = "27_EL RETIRO"
$ws.Range("D299").Value = 64
$ws.Range("A300").Value = "17:36:10"
$ws.Range("B300").Value = "19:17"
$ws.Range("C300").Value = "27_EL RETIRO"
$ws.Range("D300").Value = 101
$ws.Range("A301").Value = "17:36:10"
$ws.Range("B301").Value = "19:17"
$ws.Range("C301").Value = "14X44_ABASTO"
$ws.Range("D301").Value = 101
$ws.Range("A302").Value = "17:56:03"
$ws.Range("B302").Value = "19:21"
$ws.Range("C302").Value = "23_HERNANDEZ"
$ws.Range("D302").Value = 85
$ws.Range("A303").Value = "18:12:30"
$ws.Range("B303").Value = "19:22"
$ws.Range("C303").Value = "23_HERNANDEZ"
$ws.Range("D303").Value = 70
$ws.Range("A304").Value = "17:36:10"
$ws.Range("B304").Value = "19:28"
$ws.Range("C304").Value = "215C_EL PATO"
$ws.Range("D304").Value = 112
$ws.Range("B305").Value = "19:35"
$ws.Range("C305").Value = "11_ETCHEVERRY"
$ws.Range("D305").Value = 107
$ws.Range("B306").Value = "19:36"
$ws.Range("C306").Value = "11_ETCHEVERRY"
$ws.Range("D306").Value = 100
$ws.Range("A307").Value = "17:48:33"
$ws.Range("B307").Value = "19:39"
$ws.Range("C307").Value = "15X38_ABASTO"
$ws.Range("D307").Value = 111
$ws.Range("A308").Value = "17:56:03"
$ws.Range("B308").Value = "19:52"
$ws.Range("C308").Value = "81_EL PELIGRO"
$ws.Range("D308").Value = 116
$ws.Range("E308").Value = "LP1912"
$ws.Range("A309").Value = "17:56:03"
$ws.Range("B309").Value = "19:53"
$ws.Range("C309").Value = "225_GOMEZ"
$ws.Range("D309").Value = 117
$ws.Range("E309").Value = "LP1912"
$ws.Range("A310").Value = "18:31:25"
$ws.Range("B310").Value = "20:12"
$ws.Range("C310").Value = "215C_EL PATO"
$ws.Range("D310").Value = 101
$ws.Range("E310").Value = "LP1912"
$ws.Range("A311").Value = "18:31:25"
$ws.Range("B311").Value = "20:22"
$ws.Range("C311").Value = "15_ABASTO"
$ws.Range("D311").Value = 111
$ws.Range("E311").Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 18:31:25"
$ws.Range("A3").Value = "Total filas: 49"
$ws.Range("A54").Value = "18:31:25"
$ws.Range("B54").Value = "20:12"
$ws.Range("C54").Value = "215C_EL PATO"
$ws.Range("D54").Value = 101
$ws.Range("E54").Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 18:31:25"
